# Auto-update draw results: append the 2025-12-10 Pick 4 draw as a new
# row at the bottom of the Results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Find the next empty row right after the current data (row 84 -> row 85).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# The sheet stores every column as literal text (dates, zero-padded phase
# codes, and dash-joined results all need to stay text, not be coerced to
# numbers/dates). Format the row as Text first so Value assignment keeps
# the values as strings, then restore the default "Normal" style so the
# new row doesn't pick up a stray number-format style id.
$newRange = $ws.Range("A" + $newRow + ":E" + $newRow)
$newRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-12-10"
$ws.Cells.Item($newRow, 2).Value = "Pick 4"
$ws.Cells.Item($newRow, 3).Value = "251210"
$ws.Cells.Item($newRow, 4).Value = "9-0-6-0"
$ws.Cells.Item($newRow, 5).Value = "2025-12-10T21:45:10.720+04:00"

$newRange.Style = "Normal"
